# Add example term values to the single annotation row that will remain,
# then drop the three extra (duplicate/blank) example rows so the
# annotationTable / sheet only spans A1:O2 (header + one data row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Computation")

# Fill in example values for the "Genome Assembly Version", "data processing
# action" and "processed data file format" parameter columns on row 2.
$ws.Range("F2").Value = "B73 RefGen_v4"
$ws.Range("I2").Value = "Read count extraction and normalization were performed using CLC genomic benchwork"
$ws.Range("L2").Value = "txt"

# Remove the redundant example rows 3-5 (delete bottom-up so row indices
# of the rows still to be removed don't shift).
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(3).Delete()
